$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1059.2188
$ws.Range("J17").Value = 1082.1177
$ws.Range("L17").Value = 3246.3531
$ws.Range("N17").Value = -3582.3531
$ws.Range("H43").Value = 7217.857
$ws.Range("I43").Value = 6560.778
$ws.Range("J43").Value = 8400.6
$ws.Range("K43").Value = 6560.778
$ws.Range("L43").Value = 8400.6
$ws.Range("M43").Value = -6491.778
$ws.Range("N43").Value = -8538.6
$ws.Range("H62").Value = 2181.875
$ws.Range("I62").Value = 1867.5
$ws.Range("K62").Value = 1867.5
$ws.Range("M62").Value = -1243.5
$ws.Range("H65").Value = 2181.875
$ws.Range("I65").Value = 1867.5
$ws.Range("K65").Value = 9337.5
$ws.Range("M65").Value = -6217.5
$ws.Range("H100").Value = 292.66666
$ws.Range("I100").Value = 292.66666
$ws.Range("K100").Value = 292.66666
$ws.Range("M100").Value = 248.33334
$ws.Range("H106").Value = 31907.363
$ws.Range("I106").Value = 34298.1
$ws.Range("J106").Value = 8000
$ws.Range("K106").Value = 34298.1
$ws.Range("L106").Value = 8000
$ws.Range("M106").Value = -33667.1
$ws.Range("N106").Value = -9262
$ws.Range("H116").Value = 3750.75
$ws.Range("I116").Value = 1005
$ws.Range("K116").Value = 1005
$ws.Range("M116").Value = 2437
$ws.Range("H137").Value = 2062.625
$ws.Range("I137").Value = 1834
$ws.Range("J137").Value = 2199.8
$ws.Range("K137").Value = 5502
$ws.Range("L137").Value = 6599.400000000001
$ws.Range("M137").Value = -2952
$ws.Range("N137").Value = -11699.4
$ws.Range("H138").Value = 2622.1892
$ws.Range("I138").Value = 1966.6897
$ws.Range("K138").Value = 5900.0691
$ws.Range("M138").Value = -760.0690999999997
$ws.Range("H141").Value = 3122.111
$ws.Range("I141").Value = 1637.625
$ws.Range("K141").Value = 4912.875
$ws.Range("M141").Value = 267.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3050.985
$ws.Range("I32").Value = 2275.7378
$ws.Range("K32").Value = 2275.7378
$ws.Range("M32").Value = -1988.7378
$ws.Range("H122").Value = 1121217.5
$ws.Range("I122").Value = 1260994.6
$ws.Range("K122").Value = 3782983.8
$ws.Range("M122").Value = -3780533.8
$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H140").Value = 85500
$ws.Range("I140").Value = 65000
$ws.Range("J140").Value = 106000
$ws.Range("K140").Value = 65000
$ws.Range("L140").Value = 106000
$ws.Range("M140").Value = -59820
$ws.Range("N140").Value = -116360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3490.0625
$ws.Range("I86").Value = 3269.4546
$ws.Range("J86").Value = 3975.4
$ws.Range("K86").Value = 3269.4546
$ws.Range("L86").Value = 3975.4
$ws.Range("M86").Value = -2146.4546
$ws.Range("N86").Value = -6221.4
$ws.Range("H89").Value = 3490.0625
$ws.Range("I89").Value = 3269.4546
$ws.Range("J89").Value = 3975.4
$ws.Range("K89").Value = 16347.273
$ws.Range("L89").Value = 19877
$ws.Range("M89").Value = -10731.273
$ws.Range("N89").Value = -31109
$ws.Range("H107").Value = 2098.6667
$ws.Range("I107").Value = 2098.6667
$ws.Range("K107").Value = 2098.6667
$ws.Range("M107").Value = -178.6667000000002
$ws.Range("H134").Value = 2730.9
$ws.Range("I134").Value = 1509.0769
$ws.Range("K134").Value = 4527.2307
$ws.Range("M134").Value = -1992.2307

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H37").Value = 14026
$ws.Range("J37").Value = 14026
$ws.Range("L37").Value = 14026
$ws.Range("N37").Value = -14240

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 687.25
$ws.Range("I25").Value = 624.5
$ws.Range("K25").Value = 1873.5
$ws.Range("M25").Value = -1704.5
$ws.Range("H30").Value = 687.25
$ws.Range("I30").Value = 624.5
$ws.Range("K30").Value = 1873.5
$ws.Range("M30").Value = -1771.5
$ws.Range("H55").Value = 93663.55
$ws.Range("J55").Value = 4250
$ws.Range("L55").Value = 12750
$ws.Range("N55").Value = -13104
$ws.Range("H98").Value = 2203
$ws.Range("I98").Value = 2749
$ws.Range("K98").Value = 8247
$ws.Range("M98").Value = -6749
$ws.Range("H133").Value = 1065.6666
$ws.Range("I133").Value = 1065.6666
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 3196.9998
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = 1863.0002
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4896.909
$ws.Range("I80").Value = 3480.6667
$ws.Range("K80").Value = 3480.6667
$ws.Range("M80").Value = -2482.6667
$ws.Range("H83").Value = 4896.909
$ws.Range("I83").Value = 3480.6667
$ws.Range("K83").Value = 17403.3335
$ws.Range("M83").Value = -12411.3335
$ws.Range("H122").Value = 54331.316
$ws.Range("I122").Value = 1214.2222
$ws.Range("K122").Value = 3642.6666
$ws.Range("M122").Value = -1192.6666
$ws.Range("H126").Value = 2523.75
$ws.Range("I126").Value = 2593.2
$ws.Range("K126").Value = 7779.599999999999
$ws.Range("M126").Value = -5309.599999999999
$ws.Range("H132").Value = 3736.889
$ws.Range("I132").Value = 3736.889
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11210.667
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -8680.667000000001
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 500499.5
$ws.Range("J46").Value = 500499.5
$ws.Range("L46").Value = 500499.5
$ws.Range("N46").Value = -500875.5
$ws.Range("H82").Value = 1066.3334
$ws.Range("I82").Value = 913.6
$ws.Range("J82").Value = 1205.1818
$ws.Range("K82").Value = 913.6
$ws.Range("L82").Value = 1205.1818
$ws.Range("M82").Value = -552.6
$ws.Range("N82").Value = -1927.1818
$ws.Range("H85").Value = 1066.3334
$ws.Range("I85").Value = 913.6
$ws.Range("J85").Value = 1205.1818
$ws.Range("K85").Value = 913.6
$ws.Range("L85").Value = 1205.1818
$ws.Range("M85").Value = 334.4
$ws.Range("N85").Value = -3701.1818

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 10000
$ws.Range("J24").Value = 10000
$ws.Range("L24").Value = 10000
$ws.Range("N24").Value = -10460
$ws.Range("H96").Value = 6750
$ws.Range("J96").Value = 6750
$ws.Range("L96").Value = 6750
